$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

# Fill in the Question 3 lookup table (rows 56-61) with VLOOKUP formulas
# pulling FY17/FY18/FY19 diff values from the main table (A1:P52).
# Row 56 is entered on its own first, then 57:61 as a separate fill so the
# formulas end up in the same grouping as the source workbook.
$ws.Range("B56").Formula = "=VLOOKUP(A56,`$A`$1:`$P`$52,4)"
$ws.Range("C56").Formula = "=VLOOKUP(A56,`$A`$1:`$P`$52,9)"
$ws.Range("D56").Formula = "=VLOOKUP(A56,`$A`$1:`$P`$52,14)"

$ws.Range("B57:B61").Formula = "=VLOOKUP(A57,`$A`$1:`$P`$52,4)"
$ws.Range("C57:C61").Formula = "=VLOOKUP(A57,`$A`$1:`$P`$52,9)"
$ws.Range("D57:D61").Formula = "=VLOOKUP(A57,`$A`$1:`$P`$52,14)"

# Update the view so the visible window / selection matches the new extent
# of data that was just filled in.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("E61").Select()

$wb.Save()
